$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 3 and 4 entirely (prunes the shared strings only used by those rows)
$ws.Rows("3:4").Delete()

# Update row 2 values
$ws.Range("B2").Value = "tet"
$ws.Range("C2").Value = "test"

# D2 needs to hold the numeric-looking text "23" as a genuine text cell
# (not a number). Direct Value assignment of a digit string gets coerced
# to a number, so build it as a text formula result and paste-special the
# value only, which preserves its text type without adding a number format.
$ws.Range("Z1").Formula = "=""23"""
$ws.Range("Z1").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 23
